$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 621
$ws.Range("F10").Value = 320
$ws.Range("H10").Value = 416

$ws.Range("E11").Value = 407
$ws.Range("F11").Value = 222
$ws.Range("H11").Value = 286

$ws.Range("E12").Value = 615
$ws.Range("F12").Value = 346
$ws.Range("H12").Value = 432

$ws.Range("E13").Value = 149

$ws.Range("E15").Value = 184
$ws.Range("F15").Value = 79
$ws.Range("H15").Value = 129

$ws.Range("E16").Value = 217
$ws.Range("F16").Value = 112
$ws.Range("H16").Value = 160

$ws.Range("E17").Value = 111

$ws.Range("E21").Value = 147
$ws.Range("F21").Value = 83
$ws.Range("H21").Value = 114

$ws.Range("E22").Value = 180

$ws.Range("E23").Value = 212

$ws.Range("E24").Value = 234

$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 157
$ws.Range("H25").Value = 217

$ws.Range("E26").Value = 168
$ws.Range("F26").Value = 103
$ws.Range("H26").Value = 128

$ws.Range("E27").Value = 353
$ws.Range("F27").Value = 187
$ws.Range("H27").Value = 268

$ws.Range("E28").Value = 211
$ws.Range("F28").Value = 97
$ws.Range("H28").Value = 149

$ws.Range("E29").Value = 180
$ws.Range("F29").Value = 105
$ws.Range("H29").Value = 146

$ws.Range("E30").Value = 231
$ws.Range("F30").Value = 137
$ws.Range("H30").Value = 189

$ws.Range("E31").Value = 77

$ws.Range("E32").Value = 194

$ws.Range("E34").Value = 237

$ws.Range("E35").Value = 163

$ws.Range("E36").Value = 82

$ws.Range("E37").Value = 178

$ws.Range("E41").Value = 416

$ws.Range("E42").Value = 414
$ws.Range("F42").Value = 233
$ws.Range("H42").Value = 294

$ws.Range("E43").Value = 133

$ws.Range("E44").Value = 331
$ws.Range("F44").Value = 170
$ws.Range("H44").Value = 238

$ws.Range("E46").Value = 356
$ws.Range("F46").Value = 194
$ws.Range("H46").Value = 257

$ws.Range("E47").Value = 498
$ws.Range("F47").Value = 264
$ws.Range("H47").Value = 356

$ws.Range("E48").Value = 237
$ws.Range("F48").Value = 103
$ws.Range("H48").Value = 147

$ws.Range("E49").Value = 306
$ws.Range("F49").Value = 143
$ws.Range("H49").Value = 230

$ws.Range("E50").Value = 256
$ws.Range("F50").Value = 132
$ws.Range("H50").Value = 203

$ws.Range("E51").Value = 251
$ws.Range("F51").Value = 118
$ws.Range("H51").Value = 192
